# Table Created in XLSX.
#
# Populates "Diary No" / "Who Vs Who" / "Case Details" blocks for four
# cases. Each block occupies 10 rows starting at row 5, 15, 25 and 35:
#   row+0 -> "Diary No"   label in column A, value in merged B:C
#   row+1 -> "Who Vs Who" label in column A, value in merged B:C
#   row+2 -> "Case Details" label in column A, merged down 6 rows (A:A)
#
# The existing A5/A6/A7 cells already carry the centered style (s="1")
# used throughout this table, so every new cell's formatting is copied
# from A5 via copy / paste-special(formats) rather than rebuilt from
# scratch. Merges are created before the format paste so the paste only
# touches the merge anchor cell (matching how the first block was
# authored, where only the anchor cells have <c> entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$cases = @(
    @{ StartRow = 5;  DiaryNo = "Diary No.- 1 - 2020"; WhoVsWho = "MANISHA BISHT vs. VIKAS KANWAR" },
    @{ StartRow = 15; DiaryNo = "Diary No.- 2 - 2020"; WhoVsWho = "CITY MONTESSORI SCHOOL vs. ASHOK BHARGAVA" },
    @{ StartRow = 25; DiaryNo = "Diary No.- 1 - 2019"; WhoVsWho = "C. KALIDAS vs. THE SECRETARY TAMIL NADU PUBLIC SERVICE COMMISSION" },
    @{ StartRow = 35; DiaryNo = "Diary No.- 2 - 2019"; WhoVsWho = "P. SATISH KUMAR vs. UNION OF INDIA" }
)

foreach ($case in $cases) {
    $r = $case.StartRow

    $rowDiary = $r
    $rowWho = $r + 1
    $rowCase = $r + 2
    $rowCaseLast = $r + 7

    if ($rowDiary -ne 5) {
        # Labels in column A ("Diary No" / "Who Vs Who"); row 5/6 already
        # have these from the original workbook.
        $ws.Range("A$rowDiary").Value = "Diary No"
        $ws.Range("A5").Copy() | Out-Null
        $ws.Range("A$rowDiary").PasteSpecial($xlPasteFormats) | Out-Null

        $ws.Range("A$rowWho").Value = "Who Vs Who"
        $ws.Range("A5").Copy() | Out-Null
        $ws.Range("A$rowWho").PasteSpecial($xlPasteFormats) | Out-Null
    }

    # Values in column B, merged across B:C.
    $ws.Range("B$rowDiary").Value = $case.DiaryNo
    $ws.Range("B$rowDiary`:C$rowDiary").Merge() | Out-Null
    $ws.Range("A5").Copy() | Out-Null
    $ws.Range("B$rowDiary").PasteSpecial($xlPasteFormats) | Out-Null

    $ws.Range("B$rowWho").Value = $case.WhoVsWho
    $ws.Range("B$rowWho`:C$rowWho").Merge() | Out-Null
    $ws.Range("A5").Copy() | Out-Null
    $ws.Range("B$rowWho").PasteSpecial($xlPasteFormats) | Out-Null

    # "Case Details" label, merged down across A:A for 6 rows. Row
    # 7 (first block) already has this merge/style/value in the source
    # workbook, so leave it untouched -- re-merging an already-merged
    # range would spuriously materialize empty styled cells for the
    # merge followers (A8:A12).
    if ($rowCase -ne 7) {
        $ws.Range("A$rowCase").Value = "Case Details"
        $ws.Range("A$rowCase`:A$rowCaseLast").Merge() | Out-Null
        $ws.Range("A5").Copy() | Out-Null
        $ws.Range("A$rowCase").PasteSpecial($xlPasteFormats) | Out-Null
    }
}
